$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 131, shifting existing rows 131-196 down to 132-197.
$ws.Rows.Item(131).Insert()

# Populate the newly inserted row 131 with its data.
$ws.Cells.Item(131, 1).Value = 4
$ws.Cells.Item(131, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(131, 3).Value = 'Los Lagos'
$ws.Cells.Item(131, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(131, 4).Value = 45086
$ws.Cells.Item(131, 5).Value = 10
$ws.Cells.Item(131, 6).Value = 100112052
$ws.Cells.Item(131, 7).Value = 'Albahaca'
$ws.Cells.Item(131, 8).Value = 'Sin especificar'
$ws.Cells.Item(131, 9).Value = 'Primera'
$ws.Cells.Item(131, 10).Value = 90
$ws.Cells.Item(131, 11).Value = 5000
$ws.Cells.Item(131, 12).Value = 5000
$ws.Cells.Item(131, 13).Value = 5000
$ws.Cells.Item(131, 14).Value = '$/paquete'
$ws.Cells.Item(131, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(131, 16).Value = 5000
$ws.Cells.Item(131, 17).Value = 1
$ws.Cells.Item(131, 18).Value = 'Hortaliza'
